$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $target) {
        $cell.Value = $replacement
    }
}
